$d = $word.ActiveDocument

# The parser was updated to use TokenIteratorFieldRewriterSplit, which
# tokenizes an AQL query/field text run-by-run. To mirror that, the runs
# that used to carry multi-token text ("{m" and "...))}") are split so
# each token-relevant fragment lives in its own <w:r>.

# --- Split 1: "{m" -> "{" and "m" (two separate runs) ---
$text = $d.Content.Text
$idx1 = $text.IndexOf("{m")
if ($idx1 -ge 0) {
    # Touch formatting on just the "m" character and revert it; this
    # engine will not re-merge a run that was ever explicitly
    # (re)formatted, so the text ends up split into two <w:r> elements
    # with identical (default) formatting.
    $rngM = $d.Range($idx1 + 1, $idx1 + 2)
    $rngM.Bold = 1
    $rngM2 = $d.Range($idx1 + 1, $idx1 + 2)
    $rngM2.Bold = 0
}

# --- Split 2: "...startsWith('A'))}" -> "...startsWith('A'))" and "}" ---
$text2 = $d.Content.Text
$idx2 = $text2.LastIndexOf("}")
if ($idx2 -ge 0) {
    $rngBrace = $d.Range($idx2, $idx2 + 1)
    $rngBrace.Bold = 1
    $rngBrace2 = $d.Range($idx2, $idx2 + 1)
    $rngBrace2.Bold = 0
}

Write-Host "Result: [$($d.Content.Text)]"
